# Regenerate the "K" (strikeout) column (G) values for save_data sheet.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
#          calc and write s_vals"
# Only column G (header "K") values change, rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 6
    4  = 9
    5  = 7
    6  = 10
    7  = 9
    8  = 10
    9  = 7
    10 = 6
    11 = 6
    12 = 7
    13 = 4
    14 = 8
    15 = 9
    16 = 7
    17 = 6
    18 = 3
    19 = 6
    20 = 9
    21 = 8
    22 = 5
    23 = 5
    24 = 4
    25 = 1
    26 = 7
    27 = 4
    28 = 6
    29 = 9
    30 = 6
    31 = 5
    32 = 8
    33 = 7
    34 = 8
    35 = 7
    36 = 7
    37 = 4
    38 = 3
    39 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
